$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names (title column A)
$ws.Range("A2").Value = "Iphone 17"
$ws.Range("A3").Value = "samsung 10"

# Autofit column A to best-fit its contents
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update selection to A3
$ws.Range("A3").Select() | Out-Null
